# Started working on interrupts for wireless
#
# Fill in the newly-discovered CC2500 GDIO pin assignments on both the
# "Sender" and "Receiver" sheets, rename the old "CC2500 GDI0" label to
# "CC2500 GDIO0", and restore each sheet's scroll/selection state.

$wb = $excel.ActiveWorkbook

$sender   = $wb.Worksheets.Item("Sender")
$receiver = $wb.Worksheets.Item("Receiver")

# --- Sender sheet (Port B, pins 10 & 11) ---------------------------------
# Pin 10 (row 29) previously had no function noted; it is the second GDIO line.
$sender.Range("C29").Value = "CC2500_GDIO2"
# Pin 11 (row 30) was labelled "CC2500 GDI0"; rename to match the datasheet.
$sender.Range("C30").Value = "CC2500 GDIO0"

# --- Receiver sheet (Port B, pins 10-15) ---------------------------------
# Mirror the CC2500 SPI/GDIO wiring that was already filled in on Sender.
$receiver.Range("C29").Value = "CC2500_GDIO2"
$receiver.Range("C30").Value = "CC2500 GDIO0"
$receiver.Range("C31").Value = "CC2500_SPI_NSS"
$receiver.Range("C32").Value = "CC2500 SPI_SCK"
$receiver.Range("C33").Value = "CC2500_SPI_MISO"
$receiver.Range("C34").Value = "CC2500_SPI_MOSI"

# --- Restore view state (active cell / scroll position per sheet) -------
$sender.Activate()
$sender.Range("J19").Select()

$receiver.Activate()
$receiver.Range("K27").Select()

# Sender is the tab that was active/selected in the saved workbook.
$sender.Activate()
